$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (O1, P1) continuing the sequence 0..13 in N1,
# reusing the same header formatting (bold, centered, bordered) as the
# existing header cells.
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Fill in the new data columns for rows 2-10
$ws.Range("O2").Value = -0.8436378740411959
$ws.Range("P2").Value = -0.6202438358650769

$ws.Range("O3").Value = -0.4619268668033214
$ws.Range("P3").Value = -0.3970542775398513

$ws.Range("O4").Value = 0.03511856657370258
$ws.Range("P4").Value = -0.01127746424855717

$ws.Range("O5").Value = 0.4117156031279412
$ws.Range("P5").Value = 0.3720558220892083

$ws.Range("O6").Value = -0.3029526259961807
$ws.Range("P6").Value = -0.2941608839115258

$ws.Range("O7").Value = -0.1588957668628415
$ws.Range("P7").Value = -0.1586202220516812

$ws.Range("O8").Value = -0.4427110506103647
$ws.Range("P8").Value = -0.4349476056201385

$ws.Range("O9").Value = 0.004290132862520322
$ws.Range("P9").Value = 0.003637434527014847

$ws.Range("O10").Value = 0.008656730556103144
$ws.Range("P10").Value = 0.009177797020132241
